$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append to the time log (rows 47-54)
# Columns: A=Date, B=In, C=Out, D=Time Worked (formula), E=Description

$rowsData = @(
    @{ Row=47; Date=43517; In=0.41666666666666669; Out=0.47222222222222227; Desc="Sentiment Analysis " },
    @{ Row=48; Date=43517; In=0.67361111111111116; Out=0.70833333333333337; Desc="Meeting" },
    @{ Row=49; Date=43523; In=0.54166666666666663; Out=0.5625;              Desc="Sentiment Analysis " },
    @{ Row=50; Date=43524; In=0.4513888888888889;  Out=0.47222222222222227; Desc="Author Analysis " },
    @{ Row=51; Date=43524; In=0.69791666666666663; Out=0.72222222222222221; Desc="Author Analysis" },
    @{ Row=52; Date=43525; In=0.625;                Out=0.66666666666666663; Desc="Author Analysis" },
    @{ Row=53; Date=43526; In=0.70833333333333337; Out=0.72916666666666663; Desc="Author Analysis " },
    @{ Row=54; Date=43528; In=0.5;                  Out=0.52083333333333337; Desc="Author Analysis " }
)

# Fill in Date/In/Out/Formula columns first, row by row.
foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.In
    $ws.Cells.Item($row, 3).Value = $r.Out
    $ws.Cells.Item($row, 4).Formula = "=IF(OR(ISBLANK(B$row),ISBLANK(C$row)), `"`", C$row-B$row)"
}

# Fill Description column in the same order the strings were originally
# typed: "Author Analysis" (row 51) was entered before "Author Analysis "
# (row 50) with a trailing space, so the shared-string table ends up with
# "Author Analysis" before "Author Analysis ".
$ws.Cells.Item(47, 5).Value = "Sentiment Analysis "
$ws.Cells.Item(48, 5).Value = "Meeting"
$ws.Cells.Item(49, 5).Value = "Sentiment Analysis "
$ws.Cells.Item(51, 5).Value = "Author Analysis"
$ws.Cells.Item(50, 5).Value = "Author Analysis "
$ws.Cells.Item(52, 5).Value = "Author Analysis"
$ws.Cells.Item(53, 5).Value = "Author Analysis "
$ws.Cells.Item(54, 5).Value = "Author Analysis "

# Update sheet view to match recorded scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("E54").Select()
